$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data for the new rows 23-36 (duplicated model-split blocks with
#     real-world prediction recall/specificity numbers) ---

$ws.Range("A23").Value = "TS 5050 FS T1 1-12-2020"
$ws.Range("B23").Value = 62
$ws.Range("C23").Value = 56
$ws.Range("D23").Value = 3971
$ws.Range("E23").Value = 510
$ws.Range("F23").Value = 0.52542372881355903
$ws.Range("G23").Value = 0.88618611916982803

$ws.Range("A24").Value = "TS 5050 FS T1 1-13-2020"
$ws.Range("B24").Value = 45
$ws.Range("C24").Value = 69
$ws.Range("D24").Value = 3774
$ws.Range("E24").Value = 711
$ws.Range("F24").Value = 0.394736842105263
$ws.Range("G24").Value = 0.84147157190635402

$ws.Range("A25").Value = "TS 5050 FS T1 1-14-2020"
$ws.Range("B25").Value = 75
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 3715
$ws.Range("E25").Value = 785
$ws.Range("F25").Value = 0.75757575757575701
$ws.Range("G25").Value = 0.82555555555555504

$ws.Range("A26").Value = "TS 5050 FS T1 1-15-2020"
$ws.Range("B26").Value = 49
$ws.Range("C26").Value = 38
$ws.Range("D26").Value = 3868
$ws.Range("E26").Value = 644
$ws.Range("F26").Value = 0.56321839080459701
$ws.Range("G26").Value = 0.85726950354609899

$ws.Range("A27").Value = "TS 5050 FS T1 1-16-2020"
$ws.Range("B27").Value = 46
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 3973
$ws.Range("E27").Value = 565
$ws.Range("F27").Value = 0.75409836065573699
$ws.Range("G27").Value = 0.87549581313353897

$ws.Range("A28").Value = "TS 5050 FS T1 1-17-2020"
$ws.Range("B28").Value = 53
$ws.Range("C28").Value = 31
$ws.Range("D28").Value = 3810
$ws.Range("E28").Value = 705
$ws.Range("F28").Value = 0.63095238095238004
$ws.Range("G28").Value = 0.84385382059800595

$ws.Range("A29").Value = "TS 5050 FS T1 1-18-2020"
$ws.Range("B29").Value = 84
$ws.Range("C29").Value = 75
$ws.Range("D29").Value = 3907
$ws.Range("E29").Value = 533
$ws.Range("F29").Value = 0.52830188679245205
$ws.Range("G29").Value = 0.87995495495495402

$ws.Range("A30").Value = "TS 5050 FS T1 4-12-2020"
$ws.Range("B30").Value = 23
$ws.Range("C30").Value = 55
$ws.Range("D30").Value = 3981
$ws.Range("E30").Value = 540
$ws.Range("F30").Value = 0.29487179487179399
$ws.Range("G30").Value = 0.88055739880557404

$ws.Range("A31").Value = "TS 5050 FS T1 4-13-2020"
$ws.Range("B31").Value = 56
$ws.Range("C31").Value = 66
$ws.Range("D31").Value = 3849
$ws.Range("E31").Value = 628
$ws.Range("F31").Value = 0.45901639344262202
$ws.Range("G31").Value = 0.85972749609113197

$ws.Range("A32").Value = "TS 5050 FS T1 4-14-2020"
$ws.Range("B32").Value = 29
$ws.Range("C32").Value = 22
$ws.Range("D32").Value = 3788
$ws.Range("E32").Value = 760
$ws.Range("F32").Value = 0.56862745098039202
$ws.Range("G32").Value = 0.83289357959542598

$ws.Range("A33").Value = "TS 5050 FS T1 4-15-2020"
$ws.Range("B33").Value = 29
$ws.Range("C33").Value = 12
$ws.Range("D33").Value = 3656
$ws.Range("E33").Value = 902
$ws.Range("F33").Value = 0.707317073170731
$ws.Range("G33").Value = 0.80210618692408897

$ws.Range("A34").Value = "TS 5050 FS T1 4-16-2020"
$ws.Range("B34").Value = 40
$ws.Range("C34").Value = 27
$ws.Range("D34").Value = 3907
$ws.Range("E34").Value = 625
$ws.Range("F34").Value = 0.59701492537313405
$ws.Range("G34").Value = 0.86209179170344197

$ws.Range("A35").Value = "TS 5050 FS T1 4-17-2020"
$ws.Range("B35").Value = 39
$ws.Range("C35").Value = 15
$ws.Range("D35").Value = 3755
$ws.Range("E35").Value = 790
$ws.Range("F35").Value = 0.72222222222222199
$ws.Range("G35").Value = 0.82618261826182604

$ws.Range("A36").Value = "TS 5050 FS T1 4-18-2020"
$ws.Range("B36").Value = 21
$ws.Range("C36").Value = 33
$ws.Range("D36").Value = 3836
$ws.Range("E36").Value = 709
$ws.Range("F36").Value = 0.38888888888888801
$ws.Range("G36").Value = 0.84400440044004399

# --- Formatting: copy the style from the existing "4-5-2020..4-11-2020"
#     block (rows 6:8) which already carries font/fill/border/number-format,
#     then recolor the two new blocks (yellow, then pink) ---

$ws.Range("A6:G8").Copy() | Out-Null
$ws.Range("A23:G25").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:G8").Copy() | Out-Null
$ws.Range("A26:G29").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:G29").Interior.Color = 65535

$ws.Range("A6:G8").Copy() | Out-Null
$ws.Range("A30:G32").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:G8").Copy() | Out-Null
$ws.Range("A33:G36").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:G36").Interior.Color = 16738047

$excel.CutCopyMode = 0

# --- Re-establish the values after the paste-special round trip (formats
#     only were copied, so the literal values/strings are unaffected, but we
#     re-assert them defensively) ---

$ws.Columns("A:G").AutoFit() | Out-Null

$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("K21").Select() | Out-Null
